# tiles.xlsx - add Wall / Window / Entrance tile types and fix the
# "Bank" key's character from upper-case "E" to lower-case "e".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the existing "Bank" row's Char value (was "E", should be "e")
$ws.Range("A8").Value = "e"

# New rows describing the Wall / Entrance / Window tiles.
# Column A is written first for all three new rows, then column B,
# then column C, then column D - matching the order the values were
# originally entered in.
$ws.Range("A12").Value = "K"
$ws.Range("A14").Value = "X"
$ws.Range("A13").Value = "O"

$ws.Range("B12").Value = "Wall"
$ws.Range("B13").Value = "Window"
$ws.Range("B14").Value = "Entrance"

$ws.Range("C12").Value = "TBC"
$ws.Range("C13").Value = "TBC"
$ws.Range("C14").Value = "TBC"

$ws.Range("D12").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 0

$ws.Range("C14").Select()
